$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Trends Status" — updated species counts
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B2").Value = 98
$ws1.Range("B3").Value = 106
$ws1.Range("B4").Value = 425
$ws1.Range("C4").Value = 305
$ws1.Range("B5").Value = 185
$ws1.Range("B6").Value = 98
$ws1.Range("B7").Value = 19
$ws1.Range("B8").Value = 17

# ---------------------------------------------------------------------------
# Sheet "Range Status" — updated species counts
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Range Status")
$ws2.Range("B3").Value = 76
$ws2.Range("B4").Value = 231
$ws2.Range("B5").Value = 348

# ---------------------------------------------------------------------------
# Sheet "Priority Status" — updated species counts
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 204
$ws3.Range("B3").Value = 340
$ws3.Range("B4").Value = 404

# ---------------------------------------------------------------------------
# Sheet "Species qualification" — updated species counts
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B2").Value = 948
$ws4.Range("B3").Value = 530
$ws4.Range("B4").Value = 650
$ws4.Range("B5").Value = 948

# ---------------------------------------------------------------------------
# New sheet "SoIB-IUCN cross-tab" — added at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "SoIB-IUCN cross-tab"

$ws5.Range("A1").Value = " "
$ws5.Range("B1").Value = "High"
$ws5.Range("C1").Value = "Moderate"
$ws5.Range("D1").Value = "Low"
$ws5.Range("E1").Value = "Sum"
$ws5.Range("A1:E1").Font.Bold = $true
$ws5.Range("A1:E1").HorizontalAlignment = -4108

$ws5.Range("A2").Value = "Critically Endangered"
$ws5.Range("B2").Value = 16
$ws5.Range("C2").Value = 0
$ws5.Range("D2").Value = 0
$ws5.Range("E2").Value = 16

$ws5.Range("A3").Value = "Endangered"
$ws5.Range("B3").Value = 15
$ws5.Range("C3").Value = 0
$ws5.Range("D3").Value = 1
$ws5.Range("E3").Value = 16

$ws5.Range("A4").Value = "Vulnerable"
$ws5.Range("B4").Value = 43
$ws5.Range("C4").Value = 7
$ws5.Range("D4").Value = 2
$ws5.Range("E4").Value = 52

$ws5.Range("A5").Value = "Near Threatened"
$ws5.Range("B5").Value = 24
$ws5.Range("C5").Value = 32
$ws5.Range("D5").Value = 11
$ws5.Range("E5").Value = 67

$ws5.Range("A6").Value = "Least Concern"
$ws5.Range("B6").Value = 105
$ws5.Range("C6").Value = 297
$ws5.Range("D6").Value = 387
$ws5.Range("E6").Value = 789

$ws5.Range("A7").Value = "Not Recognised"
$ws5.Range("B7").Value = 1
$ws5.Range("C7").Value = 4
$ws5.Range("D7").Value = 3
$ws5.Range("E7").Value = 8

$ws5.Range("A8").Value = "Sum"
$ws5.Range("B8").Value = 204
$ws5.Range("C8").Value = 404
$ws5.Range("D8").Value = 340
$ws5.Range("E8").Value = 948
